# Automatische test-sync: 2025-06-17 22:09:18
# Appends the new mail-log entry (row 47) to the "Logs" sheet, extends the
# conditional formatting ranges to cover the new row, and bumps the
# "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new log entry to the Logs sheet -------------------------
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A47").Value = "Sollicitatie marketingfunctie"
$ws.Range("B47").Value = "mailmind.test@zohomail.eu"
$ws.Range("C47").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$ws.Range("D47").Value = "Overig"
$ws.Range("F47").Value = "2025-06-17 22:08:36"
$ws.Range("G47").Value = "Nee"

# --- 2. Extend the conditional formatting ranges to include row 47 ---------
$fcsD = $ws.Range("D2:D46").FormatConditions
for ($i = 1; $i -le $fcsD.Count(); $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($ws.Range("D2:D47"))
}

$fcsG = $ws.Range("G2:G46").FormatConditions
for ($i = 1; $i -le $fcsG.Count(); $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($ws.Range("G2:G47"))
}

# --- 3. Update the Dashboard summary count for "Overig" ---------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 12
